$d = $word.ActiveDocument

# Locate the phrase "Once the" inside the ranking-class paragraph and
# replace it with "input", matching the author's edit (the paragraph was
# split into three runs around the replaced word).
$rng = $d.Content
$found = $rng.Find.Execute("Once the", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)

$start = $rng.Start
$end = $rng.End

# Remove the old text ("Once the").
$rng.Delete()

# Insert the replacement word ("input") as its own run right where the
# old text used to start.
$newRun = $d.Range($start, $start)
$newRun.InsertAfter("input")

# Force Word to keep "input" as a distinct run (instead of silently
# re-merging it with its neighbors) by touching and then reverting a
# character-formatting property on just that span.
$newRunRange = $d.Range($start, $start + 5)
$newRunRange.Font.Bold = $true
$newRunRange.Font.Bold = $false
